# "implement logging global level"
# Settings sheet (sheet1) gets a bunch of new configuration rows (7-20)
# replacing the old DB_* rows and adding new logging / retry / DB settings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# --- Row 7: LogLevel -------------------------------------------------
$ws.Range("A7").Value = "LogLevel"
$ws.Range("B7").Value = "Information"
# clear the old C7 value (previously held DB_Connection's value)
$ws.Range("C7").Value = ""

# --- Row 8: LogFilePath ----------------------------------------------
$ws.Range("A8").Value = "LogFilePath"
$ws.Range("B8").Value = "Data\Logs\ProcessLog.txt"
$ws.Range("C8").Value = ""

# --- Row 9: EnableExternalLogging (boolean) ---------------------------
$ws.Range("A9").Value = "EnableExternalLogging"
$ws.Range("B9").Value = $true

# --- Row 10: ExternalLogEndpoint (hyperlink) --------------------------
$ws.Range("A10").Value = "ExternalLogEndpoint"
$ws.Hyperlinks.Add($ws.Range("B10"), "https://dc.services.visualstudio.com/v2/track")

# --- Row 11: LogTimeFormat --------------------------------------------
$ws.Range("A11").Value = "LogTimeFormat"
$ws.Range("B11").Value = "yyyy-MM-ddTHH:mm:ss.fffZ"

# --- Row 12: CorrelationIdPrefix ---------------------------------------
$ws.Range("A12").Value = "CorrelationIdPrefix"
$ws.Range("B12").Value = "CORR_"

# --- Row 13: ProcessName ------------------------------------------------
$ws.Range("A13").Value = "ProcessName"
$ws.Range("B13").Value = "CustomerProcess"

# --- Row 14: RetryCount (number) ----------------------------------------
$ws.Range("A14").Value = "RetryCount"
$ws.Range("B14").Value = 3

# --- Row 15: RetryDelaySeconds (number) ----------------------------------
$ws.Range("A15").Value = "RetryDelaySeconds"
$ws.Range("B15").Value = 5

# --- Row 16: ConnectionString (label only) --------------------------------
$ws.Range("A16").Value = "ConnectionString"

# --- Row 17: ConnectionString / connection string value -------------------
$ws.Range("A17").Value = "ConnectionString"
$ws.Range("B17").Value = "Data Source=10.20.30.82;Initial Catalog=CustomerDB;User Id=sql.user;Password=Asif@123;Encrypt=False"

# --- Row 18: DBProvider ----------------------------------------------------
$ws.Range("A18").Value = "DBProvider"
$ws.Range("B18").Value = "System.Data.SqlClient"

# --- Row 19: MaxRetryNumber (number) ----------------------------------------
$ws.Range("A19").Value = "MaxRetryNumber"
$ws.Range("B19").Value = 3

# --- Row 20: StatusToProcess -------------------------------------------------
$ws.Range("A20").Value = "StatusToProcess"
$ws.Range("B20").Value = "New"

# Update the sheet selection to match the author's final cursor position.
$ws.Range("A6:B6").Select() | Out-Null
